# Weekly update: insert the latest week's price observations for
# "Terminal La Palmera de La Serena - Zapallo" (Hortaliza / Camote) at the
# top of the data block (row 634), pushing the existing history down by two
# rows. This mirrors the "Fruta / hortaliza, semanal" commit.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two blank rows at row 634 (existing rows 634..669 shift to 636..671).
$ws.Rows.Item(634).Resize(2).Insert()

# New row 634: "1a (cosecha)" observation for the new week (date 44706).
$ws.Range("A634").Value = 8
$ws.Range("B634").Value = "Terminal La Palmera de La Serena"
$ws.Range("C634").Value = "Coquimbo"
$ws.Range("D634").Value = 44706
$ws.Range("E634").Value = 4
$ws.Range("F634").Value = 100112045
$ws.Range("G634").Value = "Zapallo"
$ws.Range("H634").Value = "Camote"
$ws.Range("I634").Value = "1a (cosecha)"
$ws.Range("J634").Value = 1800
$ws.Range("K634").Value = 700
$ws.Range("L634").Value = 750
$ws.Range("M634").Value = 725
$ws.Range("N634").Value = "`$/kilo (volumen en unidades)"
$ws.Range("O634").Value = "Región de O'Higgins"
$ws.Range("P634").Value = 725
$ws.Range("Q634").Value = 1
$ws.Range("R634").Value = "Hortaliza"

# New row 635: "2a (cosecha)" observation for the new week (date 44706).
$ws.Range("A635").Value = 8
$ws.Range("B635").Value = "Terminal La Palmera de La Serena"
$ws.Range("C635").Value = "Coquimbo"
$ws.Range("D635").Value = 44706
$ws.Range("E635").Value = 4
$ws.Range("F635").Value = 100112045
$ws.Range("G635").Value = "Zapallo"
$ws.Range("H635").Value = "Camote"
$ws.Range("I635").Value = "2a (cosecha)"
$ws.Range("J635").Value = 960
$ws.Range("K635").Value = 600
$ws.Range("L635").Value = 650
$ws.Range("M635").Value = 625
$ws.Range("N635").Value = "`$/kilo (volumen en unidades)"
$ws.Range("O635").Value = "Región de O'Higgins"
$ws.Range("P635").Value = 625
$ws.Range("Q635").Value = 1
$ws.Range("R635").Value = "Hortaliza"
